$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data per commit diff

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.447.33"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.692.47"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.04"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5530"
$ws.Range("E6").Value = "  +8.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.011"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2717"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.13"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07622"
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.561"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.681.58"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5827"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008472"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.31"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.518.68"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.970"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.99"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.262"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.77"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1309"
$ws.Range("E25").Value = "  +7.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.919"
$ws.Range("E26").Value = "  +4.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.78"
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.421"
$ws.Range("E28").Value = "  +6.80%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06338"
$ws.Range("E29").Value = "  -4.34%  "
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.596"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6244"
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.406"
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.721"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.250"
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.124.20"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01646"
$ws.Range("E40").Value = "  +3.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8842"
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.017"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.72"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.842.74"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  -5.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.58"
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.212"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05283"
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4304"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.088"
$ws.Range("E51").Value = "  +1.63%  "
